$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (D) cells keep their original text formatting (e.g. trailing zeros,
# multiple separators) instead of being auto-converted to numbers by Excel.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '27.564.18'
$ws.Range('E2').Value = '  +2.41%  '
$ws.Range('D3').Value = '1.877.00'
$ws.Range('E3').Value = '  +1.72%  '
$ws.Range('D4').Value = '1.015'
$ws.Range('E4').Value = '  +0.83%  '
$ws.Range('D5').Value = '313.76'
$ws.Range('E5').Value = '  +1.51%  '
$ws.Range('E6').Value = '  +0.91%  '
$ws.Range('D7').Value = '0.4797'
$ws.Range('E7').Value = '  +1.68%  '
$ws.Range('D8').Value = '0.3787'
$ws.Range('E8').Value = '  +2.84%  '
$ws.Range('D9').Value = '0.07386'
$ws.Range('E9').Value = '  +2.31%  '
$ws.Range('D10').Value = '0.9411'
$ws.Range('E10').Value = '  +1.74%  '
$ws.Range('D11').Value = '20.75'
$ws.Range('E11').Value = '  +5.66%  '
$ws.Range('D12').Value = '0.07874'
$ws.Range('E12').Value = '  +3.33%  '
$ws.Range('D13').Value = '1.884.09'
$ws.Range('E13').Value = '  +1.34%  '
$ws.Range('D14').Value = '5.457'
$ws.Range('E14').Value = '  +2.74%  '
$ws.Range('D15').Value = '6.605'
$ws.Range('E15').Value = '  +3.20%  '
$ws.Range('D16').Value = '91.29'
$ws.Range('E16').Value = '  +3.24%  '
$ws.Range('E17').Value = '  +0.85%  '
$ws.Range('D18').Value = '0.000008972'
$ws.Range('E18').Value = '  +3.54%  '
$ws.Range('E19').Value = '  +0.81%  '
$ws.Range('D20').Value = '14.96'
$ws.Range('E20').Value = '  +2.68%  '
$ws.Range('D21').Value = '27.596.49'
$ws.Range('E21').Value = '  +2.45%  '
$ws.Range('D22').Value = '5.146'
$ws.Range('E22').Value = '  +2.22%  '
$ws.Range('E23').Value = '  +1.18%  '
$ws.Range('E24').Value = '  +2.49%  '
$ws.Range('D25').Value = '153.75'
$ws.Range('E25').Value = '  +1.08%  '
$ws.Range('D26').Value = '18.58'
$ws.Range('E26').Value = '  +2.43%  '
$ws.Range('D27').Value = '2.023'
$ws.Range('E27').Value = '  +0.85%  '
$ws.Range('D28').Value = '116.08'
$ws.Range('E28').Value = '  +1.70%  '
$ws.Range('D29').Value = '5.022'
$ws.Range('E29').Value = '  +1.33%  '
$ws.Range('D30').Value = '0.08943'
$ws.Range('E30').Value = '  +1.20%  '
$ws.Range('D31').Value = '3.329'
$ws.Range('E31').Value = '  +0.74%  '
$ws.Range('D32').Value = '1.218'
$ws.Range('E32').Value = '  +4.31%  '
$ws.Range('D33').Value = '4.619'
$ws.Range('E33').Value = '  +2.90%  '
$ws.Range('D34').Value = '0.7520'
$ws.Range('E34').Value = '  +0.93%  '
$ws.Range('E35').Value = '  -2.64%  '
$ws.Range('D36').Value = '0.02077'
$ws.Range('E36').Value = '  +6.61%  '
$ws.Range('E37').Value = '  +3.05%  '
$ws.Range('D38').Value = '0.05312'
$ws.Range('E38').Value = '  +0.98%  '
$ws.Range('D39').Value = '3.010'
$ws.Range('D40').Value = '0.5372'
$ws.Range('E40').Value = '  +2.84%  '
$ws.Range('D41').Value = '7.116'
$ws.Range('E41').Value = '  +2.93%  '
$ws.Range('E42').Value = '  +0.93%  '
$ws.Range('D43').Value = '8.446'
$ws.Range('E43').Value = '  +2.79%  '
$ws.Range('E44').Value = '  +1.35%  '
$ws.Range('D45').Value = '0.4850'
$ws.Range('E45').Value = '  +3.28%  '
$ws.Range('D46').Value = '1.016'
$ws.Range('E46').Value = '  +0.93%  '
$ws.Range('D47').Value = '1.664'
$ws.Range('E47').Value = '  +3.72%  '
$ws.Range('D48').Value = '103.19'
$ws.Range('E48').Value = '  +1.10%  '
$ws.Range('D49').Value = '67.47'
$ws.Range('E49').Value = '  +2.98%  '
$ws.Range('D50').Value = '0.06119'
$ws.Range('E50').Value = '  +1.40%  '
$ws.Range('D51').Value = '0.9018'
$ws.Range('E51').Value = '  +1.86%  '
